$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.958.65'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.223.51'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.49'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.632'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '71.23'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.29%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.605'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.36'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +11.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0969'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.68%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '58.26'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.25'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.03%  '
$ws.Range('E14').Value = '  -0.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.552.01'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.03'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.894'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.64%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.204.07'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.797.97'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0965'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.28'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.74'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.22'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.09'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.11'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +18.88%  '
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('E29').Value = '  -1.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.85'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.72%  '
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('E33').Value = '  -1.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.55'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0739'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.70'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.08'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +14.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.08'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.93%  '
$ws.Range('E39').Value = '  +8.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.28'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.93'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '12.50'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +25.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.33'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.208'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.83'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.74'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.59%  '
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.67'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('E50').Value = '  +6.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.40'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.45%  '
